$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.781.67'
$ws.Range('E2').Value = '  -3.51%  '
$ws.Range('D3').Value = '3.274.21'
$ws.Range('E3').Value = '  -4.29%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''554.00'
$ws.Range('E5').Value = '  -4.28%  '
$ws.Range('D6').Value = '''139.84'
$ws.Range('E6').Value = '  -8.63%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.273.70'
$ws.Range('E8').Value = '  -4.28%  '
$ws.Range('D9').Value = '''0.463'
$ws.Range('E9').Value = '  -4.27%  '
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('E11').Value = '  -5.85%  '
$ws.Range('E12').Value = '  -3.86%  '
$ws.Range('D13').Value = '3.833.32'
$ws.Range('E13').Value = '  -4.31%  '
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '''26.48'
$ws.Range('E15').Value = '  -7.19%  '
$ws.Range('D16').Value = '3.270.54'
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('E17').Value = '  -5.42%  '
$ws.Range('D18').Value = '59.888.64'
$ws.Range('E18').Value = '  -3.38%  '
$ws.Range('D19').Value = '''6.04'
$ws.Range('E19').Value = '  -7.48%  '
$ws.Range('D20').Value = '''13.59'
$ws.Range('E20').Value = '  -6.42%  '
$ws.Range('D21').Value = '''8.45'
$ws.Range('E21').Value = '  -5.75%  '
$ws.Range('D22').Value = '''370.62'
$ws.Range('E22').Value = '  -2.95%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '''72.31'
$ws.Range('E24').Value = '  -3.86%  '
$ws.Range('D25').Value = '''0.527'
$ws.Range('E25').Value = '  -7.92%  '
$ws.Range('D26').Value = '3.405.28'
$ws.Range('E26').Value = '  -4.31%  '
$ws.Range('E27').Value = '  -9.56%  '
$ws.Range('E28').Value = '  -4.36%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  -8.92%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').Value = '''2.00'
$ws.Range('E32').Value = '  -5.59%  '
$ws.Range('D33').Value = '''7.39'
$ws.Range('E33').Value = '  -6.08%  '
$ws.Range('D34').Value = '''22.41'
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').Value = '''1.22'
$ws.Range('E35').Value = '  -8.57%  '
$ws.Range('D36').Value = '''166.25'
$ws.Range('E36').Value = '  -1.67%  '
$ws.Range('D37').Value = '''5.01'
$ws.Range('E37').Value = '  -9.05%  '
$ws.Range('D38').Value = '''1.51'
$ws.Range('E38').Value = '  -5.79%  '
$ws.Range('D39').Value = '''6.57'
$ws.Range('E39').Value = '  -5.73%  '
$ws.Range('D40').Value = '3.302.75'
$ws.Range('E40').Value = '  -4.43%  '
$ws.Range('D41').Value = '''0.0720'
$ws.Range('E41').Value = '  -8.35%  '
$ws.Range('D42').Value = '''25.52'
$ws.Range('E42').Value = '  -17.68%  '
$ws.Range('D43').Value = '''41.37'
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('D44').Value = '''0.742'
$ws.Range('E44').Value = '  -4.86%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '''4.08'
$ws.Range('E45').Value = '  -7.87%  '
$ws.Range('B46').Value = 'ONDO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D46').Value = '''1.11'
$ws.Range('E46').Value = '  -4.49%  '
$ws.Range('D47').Value = '''1.55'
$ws.Range('E47').Value = '  -7.47%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = '2.307.13'
$ws.Range('E49').Value = '  -9.45%  '
$ws.Range('D50').Value = '''6.29'
$ws.Range('E50').Value = '  -7.95%  '
$ws.Range('D51').Value = '''21.13'
$ws.Range('E51').Value = '  -6.06%  '
